$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# 1) Window was un-minimized
# ------------------------------------------------------------------
$win = $excel.ActiveWindow
$win.WindowState = -4143   # xlNormal

# ------------------------------------------------------------------
# 2) Column H (Status) is reformatted as Text ("@") - header + data
#    (creates the two new cellXfs entries used throughout column H)
# ------------------------------------------------------------------
$ws.Range("H1").NumberFormat = "@"
$ws.Range("H2:H53").NumberFormat = "@"

# ------------------------------------------------------------------
# 3) Row 52 (STR/USDT0000001 sell) is finalized:
#    status flips from "IN PROGRESS" to "DONE" and the remaining
#    finalize/fee/profit/duration cells are filled in.
# ------------------------------------------------------------------
$ws.Range("H52").Value = "DONE"
$ws.Range("I52").Value = 42864.821932870371
$ws.Range("J52").Value = "0.03882735 USDT (0.15%)"

$ws.Range("K52").Value = "   ~-35%"
$k52 = $ws.Range("K52").Characters(6, 3)
$k52.Font.Color = 255
$k52.Font.Name = "Calibri"
$k52.Font.Size = 11

$ws.Range("L52").Value = " 1 day"

# ------------------------------------------------------------------
# 4) New row 53: a fresh Buy of XRP (continuing the
#    XRP/USDT0000005 trade group), still "IN PROGRESS".
# ------------------------------------------------------------------
$ws.Range("A52:L52").Copy($ws.Range("A53:L53"))
$ws.Rows(53).RowHeight = 14.25

$ws.Range("A53").Value = 42865.342280092591
$ws.Range("B53").Value = "            Buy"
$ws.Range("C53").Value = "        XRP"
$ws.Range("D53").Value = "        0.184" + [char]10
$ws.Range("E53").Value = "         0.165USDT"
$ws.Range("F53").Value = "         180 XRP"
$ws.Range("G53").Value = " XRP/USDT0000005"
$ws.Range("H53").Value = "IN PROGRESS"
$ws.Range("I53").ClearContents()
$ws.Range("J53").ClearContents()
$ws.Range("K53").Value = "     "
$ws.Range("L53").ClearContents()

# ------------------------------------------------------------------
# 5) Sheet view scrolled/selected a bit differently afterwards
# ------------------------------------------------------------------
$ws.Range("C59").Select()
